$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (updates sheet tab name / workbook.xml)
$ws.Name = "Through 2022-11-03"

# Update header label in I1 (shared string)
$ws.Range("I1").Value = "2022 (through 11-03)"

# Update data values for new day's data (2022-11-11)
$ws.Range("I11").Value = 124
$ws.Range("I12").Value = 9
$ws.Range("I14").Value = 1409
